$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E18").Value = "Help Hours"
$ws.Range("E19").Value = "Help Hours"

$ws.Range("G18").Select()
